# Box Plot Updates, Color Updates Main Figures
#
# Nudges the position of the ten box-plot label textboxes (tx9..tx18) that
# live inside the unnamed group shape on slide 1. Only Left/Top move; the
# size (Width/Height) of every textbox is unchanged.
#
# Shape.Left/Shape.Top are exposed by PowerPoint as single-precision (float32)
# values measured in points (1 pt = 12700 EMU), so the literals below are
# chosen such that, after the float32 round-trip PowerPoint performs
# internally, they land back on the exact target EMU coordinate recorded in
# the underlying OOXML.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

# name -> (new Left in points, new Top in points)
$moves = @{
    "tx9"  = @(355.4631956062992,  204.41027822047243)
    "tx10" = @(383.44633521259846, 228.53500362992128)
    "tx11" = @(418.8657072913386,  260.244140488189)
    "tx12" = @(465.35011274015744, 284.3688658976378)
    "tx13" = @(503.9416198031496,  309.1594848188977)
    "tx14" = @(505.35035696062994, 336.47680665354335)
    "tx15" = @(372.8129578858268,  382.9104307007874)
    "tx16" = @(381.8197325393701,  410.22775253543307)
    "tx17" = @(320.6270448740158,  257.12413048818894)
    "tx18" = @(327.73751831496065, 284.4414523228346)
}

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if ($moves.ContainsKey($shp.Name)) {
        $xy = $moves[$shp.Name]
        $shp.Left = $xy[0]
        $shp.Top  = $xy[1]
    }
}
